$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update phone numbers in existing rows (D2, H2, H3) to the new shared value
$ws.Range("D2").Value = 251920864496
$ws.Range("H2").Value = 251920864496
$ws.Range("H3").Value = 251920864496

# Add new student/class/teacher row (row 4)
$ws.Range("A4").Value = "Samuel Ayalew"
$ws.Range("B4").Value = 19
$ws.Range("C4").Value = "M"
$ws.Range("D4").Value = 251931653440
$ws.Range("E4").Value = 12
$ws.Range("G4").Value = "Ayalew Bikago"
$ws.Range("H4").Value = 251931653440

# Update the active selection to match the edited workbook
$ws.Range("E5").Select()
